$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B/C values per the diff
$ws.Range("B3").Value = 500    # minecraft:campfire
$ws.Range("B4").Value = 500    # minecraft:fire
$ws.Range("B7").Value = 1000   # minecraft:lava
$ws.Range("C7").Value = 5
$ws.Range("B8").Value = 500    # minecraft:magma_block
$ws.Range("C8").Value = 5
$ws.Range("B9").Value = 500    # frostedheart:oil_burner
$ws.Range("B17").Value = 400   # caupona:mud_kitchen_stove
$ws.Range("B18").Value = 500   # caupona:stone_brick_kitchen_stove
$ws.Range("B19").Value = 500   # caupona:brick_kitchen_stove
$ws.Range("B20").Value = 500   # caupona:opus_incertum_kitchen_stove
$ws.Range("B21").Value = 500   # caupona:opus_latericium_kitchen_stove

# Update sheet view: zoom + selection
$ws.Select()
$excel.ActiveWindow.Zoom = 288
$ws.Range("B6").Select()
